$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = '27.912.94'
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = '  -0.21%  '
# Row 3
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = '1.853.39'
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = '  -0.81%  '
# Row 4
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = '1.005'
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = '  +0.37%  '
# Row 5
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '311.12'
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = '  -0.47%  '
# Row 6
$ws.Range("E6").Value = '  +0.20%  '
# Row 7
$ws.Range("E7").Value = '  +1.42%  '
# Row 8
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '0.3800'
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = '  -0.58%  '
# Row 9
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.08218'
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = '  -8.68%  '
# Row 10
$ws.Range("B10").Value = 'Polygon'
$ws.Range("C10").Value = 'https://coinranking.com/coin/uW2tk-ILY0ii+polygon-matic'
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '1.105'
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = '  -1.36%  '
# Row 11
$ws.Range("B11").Value = 'Polkadot'
$ws.Range("C11").Value = 'https://coinranking.com/coin/25W7FG7om+polkadot-dot'
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '6.175'
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = '  -3.44%  '
# Row 12
$ws.Range("B12").Value = 'Solana'
$ws.Range("C12").Value = 'https://coinranking.com/coin/zNZHO_Sjf+solana-sol'
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '20.40'
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = '  -1.49%  '
# Row 13
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '1.853.52'
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = '  -0.57%  '
# Row 14
$ws.Range("B14").Value = 'Chainlink'
$ws.Range("C14").Value = 'https://coinranking.com/coin/VLqpJwogdhHNb+chainlink-link'
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '7.180'
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = '  -0.82%  '
# Row 15
$ws.Range("B15").Value = 'BinanceUSD'
$ws.Range("C15").Value = 'https://coinranking.com/coin/vSo2fu9iE1s0Y+binanceusd-busd'
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '1.002'
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = '  +0.05%  '
# Row 16
$ws.Range("B16").Value = 'ShibaInu'
$ws.Range("C16").Value = 'https://coinranking.com/coin/xz24e0BjL+shibainu-shib'
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '0.00001091'
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = '  -0.99%  '
# Row 17
$ws.Range("B17").Value = 'Litecoin'
$ws.Range("C17").Value = 'https://coinranking.com/coin/D7B1x_ks7WhV5+litecoin-ltc'
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '90.19'
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = '  -0.94%  '
# Row 18
$ws.Range("B18").Value = 'TRON'
$ws.Range("C18").Value = 'https://coinranking.com/coin/qUhEFk1I61atv+tron-trx'
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '0.06601'
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = '  -1.17%  '
# Row 19
$ws.Range("B19").Value = 'Avalanche'
$ws.Range("C19").Value = 'https://coinranking.com/coin/dvUj0CzDZ+avalanche-avax'
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '17.62'
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = '  -2.37%  '
# Row 20
$ws.Range("B20").Value = 'Dai'
$ws.Range("C20").Value = 'https://coinranking.com/coin/MoTuySvg7+dai-dai'
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '1.001'
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = '  +0.19%  '
# Row 21
$ws.Range("B21").Value = 'Uniswap'
$ws.Range("C21").Value = 'https://coinranking.com/coin/_H5FVG9iW+uniswap-uni'
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '5.991'
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = '  -2.22%  '
# Row 22
$ws.Range("B22").Value = 'WrappedBTC'
$ws.Range("C22").Value = 'https://coinranking.com/coin/x4WXHge-vvFY+wrappedbtc-wbtc'
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '27.959.12'
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = '  -0.15%  '
# Row 23
$ws.Range("B23").Value = 'Cosmos'
$ws.Range("C23").Value = 'https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom'
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '10.99'
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = '  -4.70%  '
# Row 24
$ws.Range("B24").Value = 'Toncoin'
$ws.Range("C24").Value = 'https://coinranking.com/coin/67YlI0K1b+toncoin-ton'
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '2.235'
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = '  -1.99%  '
# Row 25
$ws.Range("B25").Value = 'LEO'
$ws.Range("C25").Value = 'https://coinranking.com/coin/mqtUpyBxu8O8+leo-leo'
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '3.392'
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = '  +0.08%  '
# Row 26
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '2.545'
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = '  +1.49%  '
# Row 27
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '2.068.07'
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = '  -0.88%  '
# Row 28
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '157.57'
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = '  -0.57%  '
# Row 29
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '20.38'
$ws.Range("D29").Style = "Normal"
# Row 30
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '124.54'
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = '  -1.35%  '
# Row 31
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '0.1053'
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = '  -0.74%  '
# Row 32
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '1.030'
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = '  -2.68%  '
# Row 33
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '5.577'
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = '  -0.09%  '
# Row 34
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '3.586'
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = '  +0.08%  '
# Row 35
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '9.505'
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = '  +1.11%  '
# Row 36
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '0.06519'
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = '  -0.51%  '
# Row 37
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '0.02400'
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = '  -0.90%  '
# Row 38
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '0.2157'
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = '  -1.83%  '
# Row 39
$ws.Range("E39").Value = '  -0.69%  '
# Row 40
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '1.238'
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = '  -3.71%  '
# Row 41
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '0.6380'
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = '  +0.05%  '
# Row 42
$ws.Range("B42").Value = 'Aptos'
$ws.Range("C42").Value = 'https://coinranking.com/coin/HGYj5JCv5+aptos-apt'
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '11.09'
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = '  -4.40%  '
# Row 43
$ws.Range("B43").Value = 'InternetComputer(DFINITY)'
$ws.Range("C43").Value = 'https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp'
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '4.826'
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = '  -1.53%  '
# Row 44
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '0.6040'
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = '  +0.46%  '
# Row 45
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '13.05'
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = '  -1.40%  '
# Row 46
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '1.276'
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = '  +0.09%  '
# Row 47
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '3.652'
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = '  -0.61%  '
# Row 48
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '1.976'
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = '  -1.21%  '
# Row 49
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '1.197'
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = '  -2.05%  '
# Row 50
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '119.48'
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = '  -1.43%  '
# Row 51
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '78.39'
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = '  -0.26%  '